#
# Apply hybrid bold + color ("2C3E50") highlighting to quantitative
# metrics (percentages, dollar amounts, large numbers) inside specific
# bullet / impact paragraphs, matching the target diff.

$d = $word.ActiveDocument
$bullet = [char]0x2022
$plusminus = [char]0x00B1
$highlightColor = 5258796   # COLORREF for OOXML w:color "2C3E50" (0x00503E2C)

function Find-ParagraphByPrefix($doc, $prefix) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Find-ParagraphByExactText($doc, $exact) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $exact) {
            return $p
        }
    }
    return $null
}

function Set-MetricBold($paragraph, $metricText) {
    if ($paragraph -eq $null) {
        Write-Host "WARNING: target paragraph not found, skipping metric -> $metricText"
        return
    }
    # Search only within the given paragraph's Range so occurrences of the
    # same digits elsewhere in the document are left untouched.
    $rng = $paragraph.Range
    $found = $rng.Find.Execute($metricText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: metric not found -> $metricText"
        return
    }
    $rng.Font.Bold = 1
    $rng.Font.Color = $highlightColor
}

# 1) "Discovered systematic race coding errors ... from 23% to 64%"
$prefix1 = $bullet + " Discovered systematic race coding errors"
$p1 = Find-ParagraphByPrefix $d $prefix1
Set-MetricBold $p1 "23%"
Set-MetricBold $p1 "64%"

# 2) "Achieved 87% ... industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%"
$prefix2 = $bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing"
$p2 = Find-ParagraphByPrefix $d $prefix2
Set-MetricBold $p2 "87%"
Set-MetricBold $p2 "71%"
$pm1 = $plusminus + "4.2%"
$pm2 = $plusminus + "2.1%"
Set-MetricBold $p2 $pm1
Set-MetricBold $p2 $pm2

# 3) "Wrote RFP and analyzed bids from 1,200 vendors for research platform development"
$prefix3 = $bullet + " Wrote RFP and analyzed bids from"
$p3 = Find-ParagraphByPrefix $d $prefix3
Set-MetricBold $p3 "1,200"

# 4) "Created comprehensive meta-analysis framework ... $400M ... now valued at $1B+"
$prefix4 = $bullet + " Created comprehensive meta-analysis framework"
$p4 = Find-ParagraphByPrefix $d $prefix4
Set-MetricBold $p4 '$400M'
Set-MetricBold $p4 '$1B'

# 5) "Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M"
$prefix5 = $bullet + " Algorithm reduced mapping costs by"
$p5 = Find-ParagraphByPrefix $d $prefix5
Set-MetricBold $p5 "73.5%"
Set-MetricBold $p5 '$4.7M'

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%" (short form,
#    no trailing ", reducing..." clause) - match the exact text to disambiguate from paragraph #2.
$exact6 = $bullet + " Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
$p6 = Find-ParagraphByExactText $d $exact6
Set-MetricBold $p6 "87%"
Set-MetricBold $p6 "71%"
